# Apply "added PtDemographic and new schema variables" edit.
#
# Summary of the change (search_strings sheet):
#   - Insert a new boolean "Time Series" column before the existing
#     "Pacmed ontology" column (old F -> new G, old G -> new H, old H -> new I).
#   - Insert 7 new rows after the "BMI" row (old row 4) for new
#     PtDemographics variables: Weight, Height, DOB, Ethnicity, Postcode,
#     Admit source, Hospital Number.
#   - Fix up the BMI row (now row 4) data and flag it as "possibly derived"
#     with an explanatory note.
#   - Populate the new "Time Series" column: FALSE for all PtDemographics
#     rows, TRUE for every other (already existing) variable row.
#   - Add a "Notes" header for the (shifted) Notes column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Select()

# --- Structural edits -------------------------------------------------

# Insert 7 blank rows for the new PtDemographics variables, right after
# the BMI row (old row 4 / new row 4).
$ws2.Range("A5:A11").EntireRow.Insert()

# Insert a new column for the "Time Series" flag before the old
# "Pacmed ontology" column (column F).
$ws2.Columns.Item(6).Insert()

# --- Header row ---------------------------------------------------------

$ws2.Range("F1").Value = "Time Series"
$ws2.Range("I1").Value = "Notes"

# --- Fix up the BMI row (row 4) ------------------------------------------

$ws2.Range("C4").Value = "bmi, body, index"
$ws2.Range("D4").Value = "PtDemographics"
$ws2.Range("E4").Value = 1
$ws2.Range("I4").Value = "If not recorded, use weight and height to compute (i.e. derived)"

# --- New PtDemographics rows (5-11) --------------------------------------

$ws2.Range("A5").Value = "Weight"
$ws2.Range("B5").Value = "weight"
$ws2.Range("C5").Value = "weight"
$ws2.Range("D5").Value = "PtDemographics"
$ws2.Range("E5").Value = 0

$ws2.Range("A6").Value = "Height"
$ws2.Range("B6").Value = "height"
$ws2.Range("C6").Value = "height"
$ws2.Range("D6").Value = "PtDemographics"
$ws2.Range("E6").Value = 0

$ws2.Range("A7").Value = "DOB"
$ws2.Range("B7").Value = "dob"
$ws2.Range("C7").Value = "birth, date, dob"
$ws2.Range("D7").Value = "PtDemographics"
$ws2.Range("E7").Value = 0

$ws2.Range("A8").Value = "Ethnicity"
$ws2.Range("B8").Value = "ethnicity"
$ws2.Range("C8").Value = "ethnic"
$ws2.Range("D8").Value = "PtDemographics"
$ws2.Range("E8").Value = 0

$ws2.Range("A9").Value = "Postcode"
$ws2.Range("B9").Value = "postcode"
$ws2.Range("C9").Value = "post, code, postcode, address, zip"
$ws2.Range("D9").Value = "PtDemographics"
$ws2.Range("E9").Value = 0

$ws2.Range("A10").Value = "Admit source"
$ws2.Range("B10").Value = "admit_source"
$ws2.Range("C10").Value = "admit, source, admission"
$ws2.Range("D10").Value = "PtDemographics"
$ws2.Range("E10").Value = 0

$ws2.Range("A11").Value = "Hospital Number"
$ws2.Range("B11").Value = "hospital_number"
$ws2.Range("C11").Value = "hospital, number"
$ws2.Range("D11").Value = "PtDemographics"
$ws2.Range("E11").Value = 0

# --- "Time Series" boolean column -----------------------------------------
# All of the (new) PtDemographics rows are static, not time series.
$ws2.Range("F2:F11").Value = $false
# Every other, pre-existing variable (now rows 12-39) is a time series value.
$ws2.Range("F12:F39").Value = $true
# Display booleans as TRUE/FALSE text.
$ws2.Range("F2:F39").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# --- Selection state, to mirror the authored workbook view ----------------

$ws1.Range("B20").Select()
$ws2.Select()
$ws2.Range("E11:F11").Select()
